# Beaconit.xlsx update
# - Re-label Beacon 1 placement (laiturille 4 -> laiturille 1)
# - Fill in Minor/hex/binary columns (C/F/G/H) for several Minor rows that were
#   previously left at the generic placeholder "0044" / blank, now carrying
#   concrete per-beacon Minor values with their hex/binary breakdown
# - Adjust view state (selections, scroll position, row heights) on all three sheets

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Livi" - data changes
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Livi")

# Beacon 1 placement text
$ws.Range("E4").Value = "Alikulkutunnelin pää (kulku laiturille 1)"

# Beacon 1, second iBeacon slot minor block (row 31)
$ws.Range("C31").Value = "129"
$ws.Range("F31").Value = "0000000100101001"
$ws.Range("G31").Value = "00000001"
$ws.Range("H31").Value = "00101001"

# Beacon 2, first iBeacon slot minor block (row 57)
$ws.Range("C57").Value = "22E"
$ws.Range("F57").Value = "0000001000101110"
$ws.Range("H57").Value = "00101110"

# Beacon 3, first iBeacon slot minor block (row 91)
$ws.Range("C91").Value = "22A"
$ws.Range("F91").Value = "0000001000101010"

# Beacon 4, first iBeacon slot minor block (row 115)
$ws.Range("C115").Value = "42C"
$ws.Range("F115").Value = "0000010000101100"
$ws.Range("G115").Value = "00000100"
$ws.Range("H115").Value = "00101100"

# Beacon 5, first iBeacon slot minor block (row 139)
$ws.Range("C139").Value = "430"
$ws.Range("F139").Value = "0000010000110000"
$ws.Range("G139").Value = "00000100"
$ws.Range("H139").Value = "00110000"

# Beacon 6, first iBeacon slot minor block (row 163)
$ws.Range("C163").Value = "12D"
$ws.Range("F163").Value = "0000000100101101"
$ws.Range("G163").Value = "00000001"
$ws.Range("H163").Value = "00101101"

# ---------------------------------------------------------------------------
# Sheet "HSL" - view state
# ---------------------------------------------------------------------------
$wsHSL = $wb.Worksheets.Item("HSL")
$wsHSL.Activate()
$wsHSL.Rows.Item(1).RowHeight = 64
$wsHSL.Range("A7").Select()

# ---------------------------------------------------------------------------
# Sheet "VR" - view state
# ---------------------------------------------------------------------------
$wsVR = $wb.Worksheets.Item("VR")
$wsVR.Activate()
$wsVR.Rows.Item(1).RowHeight = 64
$wsVR.Rows.Item(3).RowHeight = 48
$wsVR.Range("A5").Select()
$winVR = $excel.ActiveWindow
$winVR.ScrollRow = 3
$winVR.ScrollColumn = 1

# ---------------------------------------------------------------------------
# Sheet "Livi" - re-activate and set final view state
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("C164").Select()
$winLivi = $excel.ActiveWindow
$winLivi.ScrollRow = 141
$winLivi.ScrollColumn = 1

# ---------------------------------------------------------------------------
# Workbook window geometry (best effort; cosmetic only)
# ---------------------------------------------------------------------------
try {
    $winLivi.Left = 7140
    $winLivi.Top = 1800
    $winLivi.Width = 17420
    $winLivi.Height = 14720
} catch {}
